$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Location 1 (heading "Junior Manual QA Engineer"): the runs
#   "Junior" + <bookmarkStart _GoBack/><bookmarkEnd/> + " Manual QA "
# collapse into a single run "Junior Manual QA " and the _GoBack
# bookmark is removed from here.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$loc1 = $d.Content
$loc1.Find.Execute("Junior Manual QA ") | Out-Null
$loc1.Text = "Junior Manual QA "

# ------------------------------------------------------------------
# Location 2 (work-experience paragraph): "Junior(+) QA Manual
# Engineer ..." loses the "(" and "+)" characters (becoming
# "Junior QA Manual Engineer ...") and gains the _GoBack bookmark
# right after "Junior". The proofErr gramStart/gramEnd markers that
# wrapped "Junior(" disappear along with the edit.
# ------------------------------------------------------------------
$loc2 = $d.Content
$loc2.Find.Execute("Junior(") | Out-Null
$juniorStart = $loc2.Start
$juniorEnd = $loc2.End

# Delete the "(" first (last character of the "Junior(" match).
$parenRng = $d.Range($juniorEnd - 1, $juniorEnd)
$parenRng.Text = ""

# Delete the "+)" that immediately follows (now immediately after
# "Junior" since the "(" is gone).
$plusRng = $d.Range($juniorEnd - 1, $juniorEnd + 1)
$plusRng.Text = ""

# Insert the _GoBack bookmark right after "Junior".
$bmRng = $d.Range($juniorStart + 6, $juniorStart + 6)
$d.Bookmarks.Add("_GoBack", $bmRng)
